# Update 17-Jun-2021, end of day.
# Petty cash book: fill in transactions for 14-Jun through 17-Jun-2021
# on Sheet1 (rows 3-26), adding new "Keterangan" categories as needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 14-Jun-2021 (row 3 already dated; continues from row 3) ---
$ws.Range("D3").Formula = "=60000+260000"

$ws.Range("B4").Value = "TRANSFER BCA"
$ws.Range("D4").Formula = "=6440000+2250000+4500000+1855500+54280000+772000+12500000"

$ws.Range("B5").Value = "FREIGHT OUT"
$ws.Range("D5").Formula = "=7500+7200"

$ws.Range("B6").Value = "A/R"
$ws.Range("C6").Formula = "=64280000+1021000+672000+15045500"

$ws.Range("B7").Value = "A/P"
$ws.Range("D7").Formula = "=750000"

$ws.Range("B8").Value = "DANA KEBERSIHAN"
$ws.Range("D8").Formula = "=120000"

$ws.Range("B9").Value = "BENSIN - rush"
$ws.Range("D9").Formula = "=250000"

$ws.Range("B10").Value = "SALES - cash/retail"
$ws.Range("C10").Formula = "=5042275+19570225-600000+12500000-15045500"

$ws.Range("B11").Value = "SELISIH - lebih"
$ws.Range("C11").Value = 15500

$ws.Range("B12").Value = "SETOR KE BANK"
$ws.Range("D12").Formula = "=18000000"

# --- 15-Jun-2021 ---
$ws.Range("A13").Value = 44362
$ws.Range("B13").Value = "Wages Expense"
$ws.Range("D13").Formula = "=60000+240000"

$ws.Range("B14").Value = "TRANSFER BCA"
$ws.Range("D14").Formula = "=3150000+2200000+1814000+5500000+707000"

$ws.Range("B15").Value = "A/R"
$ws.Range("C15").Formula = "=2947500+202500+15579000"

$ws.Range("B16").Value = "SALES - cash/retail"
$ws.Range("C16").Formula = "=10456725+9555275-15579000"

$ws.Range("B17").Value = "SELISIH - lebih"
$ws.Range("C17").Value = 10000

$ws.Range("B18").Value = "SETOR KE BANK"
$ws.Range("D18").Value = 10000000

# --- 16-Jun-2021 ---
$ws.Range("A19").Value = 44363
$ws.Range("B19").Value = "Wages Expense"
$ws.Range("D19").Formula = "=60000+240000"

$ws.Range("B20").Value = "TRANSFER BCA"
$ws.Range("D20").Formula = "=11041000+1100000+1454000"

$ws.Range("B21").Value = "A/P"
$ws.Range("D21").Formula = "=700000"

$ws.Range("B22").Value = "A/R"
$ws.Range("C22").Formula = "=35333000"

$ws.Range("B23").Value = "SALES - cash/retail"
$ws.Range("C23").Formula = "=14128275+23784725-35333000"

$ws.Range("B24").Value = "SOLAR - kijang"
$ws.Range("D24").Value = 300000

$ws.Range("B25").Value = "SETOR KE BANK"
$ws.Range("D25").Value = 23000000

# --- 17-Jun-2021, end of day ---
$ws.Range("A26").Value = 44364

# Update view state: scroll the frozen pane down a bit and leave the
# cursor on the last entry of the day.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 5
$ws.Range("C26").Select()
